$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a 6th data row by copying row 5's formatting (style for G = PASS green) ---
$ws.Range("A5:G5").Copy($ws.Range("A6:G6"))

# --- Row 2: SERV_03 (Servlet: Sai mật khẩu) ---
$ws.Range("A2").Value = "SERV_03"
$ws.Range("B2").Value = "Servlet: Sai mật khẩu"
$ws.Range("C2").Value = "User: user, Pass: wrong"
$ws.Range("D2").Value = "1. Mock User CÓ tồn tại`n2. Service trả FAILED_CREDENTIALS"
$ws.Range("E2").Value = "Báo lỗi 'Mật khẩu không chính xác'"
$ws.Range("F2").Value = "OK"
$ws.Range("G2").Value = "PASS"

# --- Row 3: SERV_01 (Servlet: Đăng nhập Admin thành công) ---
$ws.Range("A3").Value = "SERV_01"
$ws.Range("B3").Value = "Servlet: Đăng nhập Admin thành công"
$ws.Range("C3").Value = "User: admin, Pass: 123"
$ws.Range("D3").Value = "1. Mock User tồn tại`n2. Mock Service trả về Admin"
$ws.Range("E3").Value = "Redirect đến Dashboard"
$ws.Range("F3").Value = "OK"
$ws.Range("G3").Value = "PASS"

# --- Row 4: SERV_02 (Servlet: Tài khoản không tồn tại) ---
$ws.Range("A4").Value = "SERV_02"
$ws.Range("B4").Value = "Servlet: Tài khoản không tồn tại"
$ws.Range("C4").Value = "User: not_exist"
$ws.Range("D4").Value = "1. Mock User KHÔNG tồn tại (DAO return false)`n2. Check lỗi"
$ws.Range("E4").Value = "Báo lỗi 'Tài khoản không tồn tại'"
$ws.Range("F4").Value = "OK"
$ws.Range("G4").Value = "PASS"

# --- Row 5: SERV_05 (User login -> Trang chủ) ---
$ws.Range("A5").Value = "SERV_05"
$ws.Range("B5").Value = "User login -> Trang chủ"
$ws.Range("C5").Value = "User: user"
$ws.Range("D5").Value = "1. Mock User tồn tại`n2. Service SUCCESS_USER"
$ws.Range("E5").Value = "Redirect /user/view-products"
$ws.Range("F5").Value = "OK"
$ws.Range("G5").Value = "PASS"

# --- Row 6: SERV_04 (Servlet: Quyền không hợp lệ) ---
$ws.Range("A6").Value = "SERV_04"
$ws.Range("B6").Value = "Servlet: Quyền không hợp lệ"
$ws.Range("C6").Value = "User: banned_user"
$ws.Range("D6").Value = "1. Mock User CÓ tồn tại`n2. Service trả FAILED_INVALID_ROLE"
$ws.Range("E6").Value = "Báo lỗi 'Quyền truy cập không hợp lệ'"
$ws.Range("F6").Value = "OK"
$ws.Range("G6").Value = "PASS"

# --- Re-run best-fit autosize on the columns whose content actually changed ---
$ws.Columns.Item(1).AutoFit()
$ws.Columns.Item(3).AutoFit()
$ws.Columns.Item(4).AutoFit()
$ws.Columns.Item(5).AutoFit()
